$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Cells whose new text would otherwise be auto-parsed as a number are
# forced to Text format first, so they stay stored as literal strings
# (matching the 'NNN.NN' style price text already used in this sheet).

$ws.Cells.Item(2, 4).Value = '60.601.61'
$ws.Cells.Item(2, 5).Value = '  +4.25%  '
$ws.Cells.Item(3, 4).Value = '2.336.65'
$ws.Cells.Item(3, 5).Value = '  +2.22%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '548.48'
$ws.Cells.Item(5, 5).Value = '  +2.53%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '132.02'
$ws.Cells.Item(6, 5).Value = '  +0.45%  '
$ws.Cells.Item(7, 5).Value = '  +0.03%  '
$ws.Cells.Item(8, 5).Value = '  -1.00%  '
$ws.Cells.Item(9, 4).Value = '2.334.42'
$ws.Cells.Item(9, 5).Value = '  +2.14%  '
$ws.Cells.Item(10, 5).Value = '  +1.56%  '
$ws.Cells.Item(11, 5).Value = '  +0.79%  '
$ws.Cells.Item(13, 5).Value = '  +1.74%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '23.86'
$ws.Cells.Item(14, 5).Value = '  +1.53%  '
$ws.Cells.Item(15, 4).Value = '2.753.33'
$ws.Cells.Item(15, 5).Value = '  +2.24%  '
$ws.Cells.Item(16, 4).Value = '60.524.36'
$ws.Cells.Item(16, 5).Value = '  +4.25%  '
$ws.Cells.Item(17, 5).Value = '  +1.28%  '
$ws.Cells.Item(18, 4).Value = '2.342.77'
$ws.Cells.Item(18, 5).Value = '  +0.61%  '
$ws.Cells.Item(19, 5).Value = '  +1.33%  '
$ws.Cells.Item(20, 5).Value = '  -0.37%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '315.66'
$ws.Cells.Item(21, 5).Value = '  +0.83%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.67'
$ws.Cells.Item(22, 5).Value = '  +3.36%  '
$ws.Cells.Item(23, 5).Value = '  -0.34%  '
$ws.Cells.Item(24, 5).Value = '  +1.87%  '
$ws.Cells.Item(25, 5).Value = '  +1.22%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 5).Value = '  -3.85%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '7.87'
$ws.Cells.Item(27, 5).Value = '  -1.15%  '
$ws.Cells.Item(28, 5).Value = '  +8.02%  '
$ws.Cells.Item(29, 5).Value = '  +12.35%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '173.44'
$ws.Cells.Item(30, 5).Value = '  +1.62%  '
$ws.Cells.Item(31, 5).Value = '  +2.66%  '
$ws.Cells.Item(32, 4).Value = '0.0₃0739'
$ws.Cells.Item(32, 5).Value = '  +2.37%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.96'
$ws.Cells.Item(33, 5).Value = '  +3.44%  '
$ws.Cells.Item(34, 5).Value = '  +11.48%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.381'
$ws.Cells.Item(35, 5).Value = '  +0.26%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '17.95'
$ws.Cells.Item(36, 5).Value = '  +0.58%  '
$ws.Cells.Item(37, 5).Value = '  +0.00%  '
$ws.Cells.Item(38, 5).Value = '  -0.07%  '
$ws.Cells.Item(39, 5).Value = '  +5.19%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '327.56'
$ws.Cells.Item(40, 5).Value = '  +13.79%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.55'
$ws.Cells.Item(41, 5).Value = '  +3.54%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '38.13'
$ws.Cells.Item(42, 5).Value = '  -0.60%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '140.15'
$ws.Cells.Item(43, 5).Value = '  +0.18%  '
$ws.Cells.Item(44, 5).Value = '  +1.50%  '
$ws.Cells.Item(45, 5).Value = '  -0.66%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '19.40'
$ws.Cells.Item(46, 5).Value = '  +7.28%  '
$ws.Cells.Item(47, 5).Value = '  +0.97%  '
$ws.Cells.Item(48, 5).Value = '  +1.88%  '
$ws.Cells.Item(49, 4).Value = '0.0₆0225'
$ws.Cells.Item(49, 5).Value = '  +22.66%  '
$ws.Cells.Item(50, 5).Value = '  +1.78%  '
$ws.Cells.Item(51, 5).Value = '  +0.72%  '
